$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1297.5
$ws.Range("I29").Value = 1297.5
$ws.Range("K29").Value = 3892.5
$ws.Range("M29").Value = -3611.5

$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2064
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -10320
$ws.Range("N77").ClearContents()

$ws.Range("H107").Value = 1765.6471
$ws.Range("I107").Value = 1901.5
$ws.Range("J107").Value = 1131.6666
$ws.Range("K107").Value = 1901.5
$ws.Range("L107").Value = 1131.6666
$ws.Range("M107").Value = 18.5
$ws.Range("N107").Value = -4971.6666

$ws.Range("H115").Value = 2182.2222
$ws.Range("I115").Value = 177.14285
$ws.Range("K115").Value = 531.4285500000001
$ws.Range("M115").Value = 1035.57145

$ws.Range("H127").Value = 1862.25
$ws.Range("I127").Value = 1268.3334
$ws.Range("K127").Value = 3805.0002
$ws.Range("M127").Value = 1154.9998

$ws.Range("H129").Value = 2068.8
$ws.Range("I129").Value = 1287.3334
$ws.Range("J129").Value = 2403.7144
$ws.Range("K129").Value = 3862.0002
$ws.Range("L129").Value = 7211.1432
$ws.Range("M129").Value = 1137.9998
$ws.Range("N129").Value = -17211.1432

$ws.Range("H135").Value = 516.0769
$ws.Range("I135").Value = 516.0769
$ws.Range("K135").Value = 4644.6921
$ws.Range("M135").Value = -2109.6921


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3773.3408
$ws.Range("J45").Value = 2033.6
$ws.Range("L45").Value = 2033.6
$ws.Range("N45").Value = -2787.6

$ws.Range("H74").Value = 2616.9138
$ws.Range("I74").Value = 2185.6487
$ws.Range("K74").Value = 2185.6487
$ws.Range("M74").Value = -1311.6487

$ws.Range("H77").Value = 2616.9138
$ws.Range("I77").Value = 2185.6487
$ws.Range("K77").Value = 10928.2435
$ws.Range("M77").Value = -6560.2435

$ws.Range("H81").Value = 1000000000
$ws.Range("J81").Value = 1000000000
$ws.Range("L81").Value = 1000000000
$ws.Range("N81").Value = -1000001996

$ws.Range("H84").Value = 1000000000
$ws.Range("J84").Value = 1000000000
$ws.Range("L84").Value = 3000000000
$ws.Range("N84").Value = -3000009984

$ws.Range("H132").Value = 2545
$ws.Range("I132").Value = 2545
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7635
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5105
$ws.Range("N132").ClearContents()

$ws.Range("H135").Value = 221864.47
$ws.Range("J135").Value = 221864.47
$ws.Range("L135").Value = 221864.47
$ws.Range("N135").Value = -232004.47


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H86").Value = 659703.94
$ws.Range("I86").Value = 5400.6
$ws.Range("K86").Value = 5400.6
$ws.Range("M86").Value = -4277.6

$ws.Range("H89").Value = 659703.94
$ws.Range("I89").Value = 5400.6
$ws.Range("K89").Value = 27003
$ws.Range("M89").Value = -21387

$ws.Range("H94").Value = 3168
$ws.Range("I94").Value = 2585.111
$ws.Range("K94").Value = 2585.111
$ws.Range("M94").Value = -2134.111

$ws.Range("H107").Value = 2926.6863
$ws.Range("I107").Value = 1805.9318
$ws.Range("J107").Value = 9971.429
$ws.Range("K107").Value = 1805.9318
$ws.Range("L107").Value = 9971.429
$ws.Range("M107").Value = 114.0681999999999
$ws.Range("N107").Value = -13811.429


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1589.8889
$ws.Range("I16").Value = 1374.6
$ws.Range("K16").Value = 1374.6
$ws.Range("M16").Value = -1087.6

$ws.Range("H31").Value = 3023.8
$ws.Range("J31").Value = 5739.5
$ws.Range("L31").Value = 5739.5
$ws.Range("N31").Value = -6329.5

$ws.Range("H34").Value = 3023.8
$ws.Range("J34").Value = 5739.5
$ws.Range("L34").Value = 5739.5
$ws.Range("N34").Value = -6143.5

$ws.Range("H113").Value = 1589.8889
$ws.Range("I113").Value = 1374.6
$ws.Range("K113").Value = 1374.6
$ws.Range("M113").Value = 795.4000000000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 214.66667
$ws.Range("I44").Value = 260.5
$ws.Range("J44").Value = 123
$ws.Range("K44").Value = 781.5
$ws.Range("L44").Value = 369
$ws.Range("M44").Value = -383.5
$ws.Range("N44").Value = -1165

$ws.Range("H57").Value = 4058.9092
$ws.Range("I57").Value = 599.75
$ws.Range("J57").Value = 6035.5713
$ws.Range("K57").Value = 1799.25
$ws.Range("L57").Value = 18106.7139
$ws.Range("M57").Value = -1240.25
$ws.Range("N57").Value = -19224.7139

$ws.Range("H82").Value = 58663.5
$ws.Range("I82").Value = 57000
$ws.Range("K82").Value = 171000
$ws.Range("M82").Value = -170594

$ws.Range("H85").Value = 58663.5
$ws.Range("I85").Value = 57000
$ws.Range("K85").Value = 171000
$ws.Range("M85").Value = -169596

$ws.Range("H99").Value = 12923.134
$ws.Range("I99").Value = 5534.9
$ws.Range("J99").Value = 27699.6
$ws.Range("K99").Value = 16604.7
$ws.Range("L99").Value = 83098.79999999999
$ws.Range("M99").Value = -14358.7
$ws.Range("N99").Value = -87590.79999999999

$ws.Range("H119").Value = 1370.4286
$ws.Range("I119").Value = 527.9
$ws.Range("K119").Value = 1583.7
$ws.Range("M119").Value = 3254.3

$ws.Range("H122").Value = 2386.875
$ws.Range("I122").Value = 637
$ws.Range("J122").Value = 2636.8572
$ws.Range("K122").Value = 5733
$ws.Range("L122").Value = 23731.7148
$ws.Range("M122").Value = -3283
$ws.Range("N122").Value = -28631.7148

$ws.Range("H123").Value = 997
$ws.Range("I123").Value = 997
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 2991
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -541
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 2789.3572
$ws.Range("J132").Value = 3013.6365
$ws.Range("L132").Value = 27122.7285
$ws.Range("N132").Value = -32182.7285


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 7822.3076
$ws.Range("J2").Value = 265.6
$ws.Range("L2").Value = 265.6
$ws.Range("N2").Value = -491.6

$ws.Range("H19").Value = 3900
$ws.Range("I19").Value = 5000
$ws.Range("K19").Value = 5000
$ws.Range("M19").Value = -4712

$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 998
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 998
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -469
$ws.Range("N22").Value = -2058

$ws.Range("H113").Value = 7942.5557
$ws.Range("I113").Value = 7605.591
$ws.Range("K113").Value = 7605.591
$ws.Range("M113").Value = -5435.591

$ws.Range("H122").Value = 4301.5625
$ws.Range("I122").Value = 3039.4736
$ws.Range("K122").Value = 9118.4208
$ws.Range("M122").Value = -6668.4208

$ws.Range("H132").Value = 1481.4286
$ws.Range("I132").Value = 1474.64
$ws.Range("J132").Value = 1538
$ws.Range("K132").Value = 4423.92
$ws.Range("L132").Value = 4614
$ws.Range("M132").Value = -1893.92
$ws.Range("N132").Value = -9674

$ws.Range("H135").Value = 121389.5
$ws.Range("J135").Value = 121389.5
$ws.Range("L135").Value = 121389.5
$ws.Range("N135").Value = -131529.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1623.6154
$ws.Range("I16").Value = 1623.6154
$ws.Range("K16").Value = 1623.6154
$ws.Range("M16").Value = -1453.6154

$ws.Range("H40").Value = 4876.654
$ws.Range("I40").Value = 3711
$ws.Range("J40").Value = 7499.375
$ws.Range("K40").Value = 3711
$ws.Range("L40").Value = 7499.375
$ws.Range("M40").Value = -3575
$ws.Range("N40").Value = -7771.375

$ws.Range("H61").Value = 7833.3335
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4798

$ws.Range("H82").Value = 1970.069
$ws.Range("I82").Value = 677.7857
$ws.Range("J82").Value = 3176.2
$ws.Range("K82").Value = 677.7857
$ws.Range("L82").Value = 3176.2
$ws.Range("M82").Value = -316.7857
$ws.Range("N82").Value = -3898.2

$ws.Range("H85").Value = 1970.069
$ws.Range("I85").Value = 677.7857
$ws.Range("J85").Value = 3176.2
$ws.Range("K85").Value = 677.7857
$ws.Range("L85").Value = 3176.2
$ws.Range("M85").Value = 570.2143
$ws.Range("N85").Value = -5672.2

$ws.Range("H113").Value = 7833.3335
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830

$ws.Range("H132").Value = 3070.5356
$ws.Range("I132").Value = 2432.389
$ws.Range("K132").Value = 7297.167
$ws.Range("M132").Value = -4767.167


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 11396.929
$ws.Range("I38").Value = 10414.272
$ws.Range("K38").Value = 10414.272
$ws.Range("M38").Value = -9941.272

$ws.Range("H100").Value = 803.44446
$ws.Range("I100").Value = 537.6
$ws.Range("K100").Value = 1075.2
$ws.Range("M100").Value = -534.2

$ws.Range("H113").Value = 953.32
$ws.Range("I113").Value = 803.6111
$ws.Range("J113").Value = 1338.2858
$ws.Range("K113").Value = 2410.8333
$ws.Range("L113").Value = 4014.8574
$ws.Range("M113").Value = -240.8332999999998
$ws.Range("N113").Value = -8354.8574

$ws.Range("H136").Value = 3516.65
$ws.Range("I136").Value = 817.2727
$ws.Range("K136").Value = 2451.8181
$ws.Range("M136").Value = 98.18190000000004

